# feat: Update Issue 3 (Product dan Sales, Daily Audit Not Finished Yet)
#
# The only content-level change in this revision is renaming the single
# worksheet from "Report Growth" to "Report Sales Details" so the template
# name matches its actual purpose (sales detail report instead of growth
# report). The rest of the source diff (fileVersion/rupBuild, the
# xr:revisionPtr session GUID, and the bookViews window position/size)
# are Excel-session/machine artifacts rewritten automatically whenever the
# authoring machine saves the file and are not meaningful content edits.

$wb = $excel.ActiveWorkbook

# Rename the (only/active) worksheet.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Report Sales Details"
